# CCC19 Derived Variables Spreadsheet - edit script
# 1) Fill in row 44 (Ca04a2 / cd20_3mo) with the Category, Description and Values
#    that were previously blank.
# 2) Insert a new row at 159 for the new "Comp41 / coinfection_other" variable
#    ("Other co-infection within +/- 2 weeks of COVID-19 dx"), pushing the
#    remainder of the table down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Complete row 44 ---------------------------------------------------
$ws.Range("C44").Value = "Cancer treatment"
$ws.Range("D44").Value = "anti-CD20 antibody within 3 months (does not necessarily catch masked regimens)"
$ws.Range("E44").Value = "0 = No; 1 = Yes"

# --- 2. Insert the new "coinfection_other" row at 159 ---------------------
$ws.Rows.Item(159).Insert()

$ws.Range("A159").Value = "Comp41"
$ws.Range("B159").Value = "coinfection_other"
$ws.Range("C159").Value = "Complications"
$ws.Range("D159").Value = "Other co-infection within +/- 2 weeks of COVID-19 dx"
$ws.Range("E159").Value = "0 = No; 1 = Yes; 99 = Unknown"

# --- Keep the backing table / autofilter in sync with the extra row -------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E324"))

# --- Mirror the saved view state (scroll position + active selection) ----
$excel.ActiveWindow.ScrollRow = 146
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("E159").Select() | Out-Null
